# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets, flips the
# Status from "Ready for handoff" to "Handed back: in sync with en-US",
# adds hyperlinks on the newly-populated target-file cells, and widens a
# few columns so the longer strings are readable.

$wb = $excel.ActiveWorkbook

$mdFileName = "03e11c4c-a08f-4387-9923-82bd72a53dbf.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/4f0edb507f4ad80c18368c002534df682c74d621/e2e/03e11c4c-a08f-4387-9923-82bd72a53dbf.md"
$statusHandedBack = "Handed back: in sync with en-US"

# ColumnWidth is quantized by the host to 1/6-character steps before the
# 5/6 "standard padding" offset is stored in the sheet XML, so feed it
# values expressed in those 1/6 steps to land exactly on the intended
# stored widths (29.9777... -> 30, 40 -> 40).
$wideColWidth = 29.166666666666668
$fullColWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F) and
# reflect the new status there too (E2/F2 mirror the per-locale Status)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = $wideColWidth   # Status
$zhcn.Columns.Item(9).ColumnWidth = $fullColWidth   # Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = $fullColWidth  # Latest Handback File

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("J2").Value = "03e11c4c-a08f-4387-9923-82bd72a53dbf.dc0d75ec29b062cf252d2f097c0ae1fa7f8445ca.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-13 07:14:39"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = $wideColWidth   # Status
$dede.Columns.Item(9).ColumnWidth = $fullColWidth   # Latest Target File
$dede.Columns.Item(10).ColumnWidth = $fullColWidth  # Latest Handback File

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("J2").Value = "03e11c4c-a08f-4387-9923-82bd72a53dbf.dc0d75ec29b062cf252d2f097c0ae1fa7f8445ca.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 07:14:49"

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
